# Write three text values into column D of the active sheet, then
# select the last cell written (matching the workbook's recorded
# selection state).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "asldfkj"
$ws.Range("D10").Value = "saldjf"
$ws.Range("D13").Value = "asdlkfj"

$ws.Range("D13").Select()
